$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 updated values
$ws.Range("B7").Value = 0.07560323850909009
$ws.Range("C7").Value = 0.2548347997162139
$ws.Range("D7").Value = 0.1052628688958802
$ws.Range("E7").Value = 0.3244423968840697
$ws.Range("F7").Value = 0.3346496625213748
$ws.Range("G7").Value = 9

# Row 8 updated values
$ws.Range("B8").Value = 0.1266874699694346
$ws.Range("C8").Value = 0.3332451476001554
$ws.Range("D8").Value = 0.1748387843367057
$ws.Range("E8").Value = 0.418137279295575
$ws.Range("F8").Value = 0.4226555369927508
$ws.Range("G8").Value = 9
